$wb = $excel.ActiveWorkbook

# --- Position 1: -> RS_TO-1651589027552185 (2 rows) ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "RS_TO-1651589027552185"
$ws.Range("A4:B5").Clear()
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eyes closed"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eyes open"

# --- Position 2: -> GNG_TO-16515890275990286 (4 rows) ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "GNG_TO-16515890275990286"
$ws.Range("A6:B10").Clear()
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "go_stims-16515890275678012.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "GNG_stims-1651589027583404.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "go_stims-1651589027583404.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "GNG_stims-16515890275990286.csv"

# --- Position 3: -> vSAT_TO-16515890276615295 (4 rows) ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "vSAT_TO-16515890276615295"
$ws.Range("A2").Copy($ws.Range("A4:A5"))
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "vSAT_stims-16515890276302774.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "SAT_stims-16515890276146667.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "SAT_stims-16515890275990286.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "vSAT_stims-1651589027645902.csv"

# --- Position 4: -> NB_TO-16515890296960516 (9 rows) ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "NB_TO-16515890296960516"
$ws.Range("A2").Copy($ws.Range("A8:A10"))
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "TB-16515890293805003.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ZB-match_1-16515890281873164.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "ZB-match_7-1651589027810048.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "TB-16515890292836802.csv"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "TB-16515890296804595.csv"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "OB-16515890288685343.csv"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "ZB-match_1-16515890285538242.csv"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "OB-16515890285850759.csv"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "OB-16515890289341435.csv"

# --- Position 5: -> TOL_TO-1651589029742927 (6 rows) ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "TOL_TO-1651589029742927"
$ws.Range("A2").Copy($ws.Range("A6:A7"))
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "MM_stims-1651589029711677.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ZM_stims-16515890296960516.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "MM_stims-16515890297273018.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "ZM_stims-1651589029711677.csv"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "MM_stims-1651589029742927.csv"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "ZM_stims-16515890297273018.csv"
